$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ------------------------------------------------------------------
# Phase 1: set up formatting for the two new rows (22 and 23) by
# cloning formats from existing rows that already carry the desired
# look (row 21 as a base template, then row 4 / G4 to fix up a few
# columns that differ from row 21's template).
# ------------------------------------------------------------------
$ws.Range("A21:K21").Copy() | Out-Null
$ws.Range("A22:K22").PasteSpecial(-4122) | Out-Null
$ws.Range("A21:K21").Copy() | Out-Null
$ws.Range("A23:K23").PasteSpecial(-4122) | Out-Null

$ws.Range("A4:D4").Copy() | Out-Null
$ws.Range("A22:D22").PasteSpecial(-4122) | Out-Null
$ws.Range("A4:D4").Copy() | Out-Null
$ws.Range("A23:D23").PasteSpecial(-4122) | Out-Null

$ws.Range("G4").Copy() | Out-Null
$ws.Range("G22").PasteSpecial(-4122) | Out-Null
$ws.Range("G4").Copy() | Out-Null
$ws.Range("G23").PasteSpecial(-4122) | Out-Null

$ws.Application.CutCopyMode = $false

$ws.Rows("22").RowHeight = 120
$ws.Rows("23").RowHeight = 120

# ------------------------------------------------------------------
# Phase 2: bring over the repeated content (these all reuse existing
# shared strings, so copy the values across instead of re-typing the
# text to guarantee an exact, byte-for-byte match).
# ------------------------------------------------------------------
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A22").PasteSpecial(-4163) | Out-Null
$ws.Range("A4").Copy() | Out-Null
$ws.Range("A23").PasteSpecial(-4163) | Out-Null

$ws.Range("C4").Copy() | Out-Null
$ws.Range("C22").PasteSpecial(-4163) | Out-Null
$ws.Range("C4").Copy() | Out-Null
$ws.Range("C23").PasteSpecial(-4163) | Out-Null

$ws.Range("D4").Copy() | Out-Null
$ws.Range("D22").PasteSpecial(-4163) | Out-Null
$ws.Range("D4").Copy() | Out-Null
$ws.Range("D23").PasteSpecial(-4163) | Out-Null

$ws.Range("J4").Copy() | Out-Null
$ws.Range("J22").PasteSpecial(-4163) | Out-Null

$ws.Range("K4").Copy() | Out-Null
$ws.Range("K22").PasteSpecial(-4163) | Out-Null

$ws.Application.CutCopyMode = $false

# ID_Dato numbers
$ws.Range("B22").Value = 21
$ws.Range("B23").Value = 22

# Fecha consulta (shared text cell, edited in place so it updates the
# existing shared string used by both H22 and H23)
$ws.Range("H22").Value = "13/06/2020"
$ws.Range("H23").Value = "13/06/2020"
$ws.Range("G4").Copy() | Out-Null
$ws.Range("H22").PasteSpecial(-4122) | Out-Null
$ws.Range("G4").Copy() | Out-Null
$ws.Range("H23").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# Fecha publicacion (real dates)
$ws.Range("I22").Value = "12/6/2020"
$ws.Range("I23").Value = "11/6/2020"

# New description text for row 23
$ws.Range("F23").Value = "El trámite para el registro y emisión de los salvoconductos específicamente para la circulación de los trabajadores domésticos durante la cuarentena, decretada por las autoridades sanitarias en la ciudad de Panamá y Panamá Oeste se realizará a través de la línea telefónica de atención ciudadana del Ministerio de Trabajo y Desarrollo Laboral (Mitradel)."

# ------------------------------------------------------------------
# New hyperlinks: G22, G23 and E23. Hyperlinks.Add resets the cell
# style to the built-in Hyperlink style, so restore the intended
# format immediately afterwards.
# ------------------------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("G22"), "https://www.mitradel.gob.pa") | Out-Null
$ws.Hyperlinks.Add($ws.Range("G23"), "https://www.mitradel.gob.pa") | Out-Null
$ws.Hyperlinks.Add($ws.Range("E23"), "https://www.mitradel.gob.pa/solicitud-de-salvoconducto-para-trabajadores-domesticos-sera-via-telefonica/") | Out-Null

$ws.Range("G4").Copy() | Out-Null
$ws.Range("G22").PasteSpecial(-4122) | Out-Null
$ws.Range("G4").Copy() | Out-Null
$ws.Range("G23").PasteSpecial(-4122) | Out-Null
$ws.Range("E21").Copy() | Out-Null
$ws.Range("E23").PasteSpecial(-4122) | Out-Null
$ws.Application.CutCopyMode = $false

# ------------------------------------------------------------------
# Other edits from the diff: Fecha publicacion on row 21 changes, and
# the worksheet table needs to grow to include the new row.
# ------------------------------------------------------------------
$ws.Range("I21").Value = "7/6/2020"

$lo = $ws.ListObjects.Item(1)
$lo.Resize($ws.Range("A1:K23")) | Out-Null

$ws.Range("C2:C23").Select() | Out-Null

$ws.Range("L23").Select() | Out-Null
